# Update the "取得日時" (retrieved-at) timestamp column on the "ランサーズ"
# sheet for every existing data row (rows 2-15), reflecting a fresh
# fetch/append performed at 2025-10-15 18:34:00.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-15 18:34:00"

$ws.Range("A2:A15").Value = $newTimestamp
